# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: column letter -> new value
$values = @{
    "A" = 45921
    "B" = 85.8
    "C" = 82.84999999999999
    "D" = 75.94
    "E" = 79
    "F" = 78.97
    "G" = 75.34999999999999
    "H" = 75.94
    "I" = 85
    "J" = 69.38
    "K" = 26.2
    "L" = 10
    "M" = 0.65
    "N" = 0
    "O" = 0
    "P" = -0.01
    "Q" = -0.5
    "R" = -0.62
    "S" = -0.01
    "T" = 3.52
    "U" = 43.16
    "V" = 69.95999999999999
    "W" = 61.4
    "X" = 54.8
    "Y" = 38.09
    "Z" = 42.29
    "AA" = "0h-4h"
    "AB" = 80.90000000000001
    "AC" = "0h-2h"
    "AD" = 84.31999999999999
    "AE" = "6h-8h"
    "AF" = 80.47
    "AG" = "9h-23h"
}

foreach ($col in $values.Keys) {
    $ws.Range("$col`2").Value = $values[$col]
}
